$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D3: mayur.charvande+2@fxbytes.com -> mayur.charvande+3@fxbytes.com
$ws.Range("D3").Value = "mayur.charvande+3@fxbytes.com"

# D2: mayur.charvande+1@fxbytes.com -> mayur.charvande+4@fxbytes.com
$ws.Range("D2").Value = "mayur.charvande+4@fxbytes.com"

# Column D width change (closest reachable value to the recorded 34.28515625)
$ws.Columns.Item(4).ColumnWidth = 33.5

# Update the active selection to match
$ws.Range("D9").Select()
